$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing column C (adjusted-close recalculation) for rows 2-45 ---
$ws.Range("C2").Value = 21.72413063049316
$ws.Range("C3").Value = 21.86156272888184
$ws.Range("C4").Value = 22.27386093139648
$ws.Range("C5").Value = 22.33275985717773
$ws.Range("C6").Value = 21.79284858703613
$ws.Range("C7").Value = 21.83211326599121
$ws.Range("C8").Value = 21.82229614257812
$ws.Range("C9").Value = 21.53761672973633
$ws.Range("C10").Value = 21.39036750793457
$ws.Range("C11").Value = 20.96825408935547
$ws.Range("C12").Value = 20.9486198425293
$ws.Range("C13").Value = 20.53632164001465
$ws.Range("C14").Value = 20.27127456665039
$ws.Range("C15").Value = 19.98659324645996
$ws.Range("C16").Value = 19.82952880859375
$ws.Range("C17").Value = 19.73136138916016
$ws.Range("C18").Value = 19.65282821655273
$ws.Range("C19").Value = 19.1619987487793
$ws.Range("C20").Value = 19.43686294555664
$ws.Range("C21").Value = 19.78044319152832
$ws.Range("C22").Value = 19.87861061096191
$ws.Range("C23").Value = 19.79026031494141
$ws.Range("C24").Value = 19.82952880859375
$ws.Range("C25").Value = 19.81971168518066
$ws.Range("C26").Value = 19.81971168518066
$ws.Range("C27").Value = 19.43686294555664
$ws.Range("C28").Value = 19.64301300048828
$ws.Range("C29").Value = 19.5350284576416
$ws.Range("C30").Value = 19.34851264953613
$ws.Range("C31").Value = 19.15218162536621
$ws.Range("C32").Value = 19.26016616821289
$ws.Range("C33").Value = 18.69080352783203
$ws.Range("C34").Value = 18.44538879394531
$ws.Range("C35").Value = 18.16070556640625
$ws.Range("C36").Value = 18.08217239379883
$ws.Range("C37").Value = 18.45520210266113
$ws.Range("C38").Value = 18.72025108337402
$ws.Range("C39").Value = 18.12143898010254
$ws.Range("C40").Value = 18.55336952209473
$ws.Range("C41").Value = 18.08217239379883
$ws.Range("C42").Value = 18.11162376403809
$ws.Range("C43").Value = 19.10309791564941
$ws.Range("C44").Value = 18.87731742858887
$ws.Range("C45").Value = 19.56447982788086

# --- Append two new rows of freshly-downloaded quotes (rows 46 and 47) ---
$ws.Range("A46").Value = 44774
$ws.Range("B46").Value = 68.08000183105469
$ws.Range("C46").Value = 19.40741348266602
$ws.Range("D46").Value = 13.01000022888184
$ws.Range("E46").Value = 11.5
$ws.Range("F46").Value = 32.40000152587891
$ws.Range("G46").Value = 9.159999847412109
$ws.Range("H46").Value = 4.889999866485596
$ws.Range("I46").Value = 17.13999938964844
$ws.Range("J46").Value = 13.4399995803833
$ws.Range("K46").Value = 22.69000053405762
$ws.Range("L46").Value = 2.720000028610229
$ws.Range("M46").Value = 35.68000030517578

$ws.Range("A47").Value = 44775
$ws.Range("B47").Value = 69.15000152587891
# Column C is intentionally left blank for row 47 (no data point was returned for this date)
$ws.Range("D47").Value = 13.01000022888184
$ws.Range("E47").Value = 11.55000019073486
$ws.Range("F47").Value = 32.40000152587891
$ws.Range("G47").Value = 9.199999809265137
$ws.Range("H47").Value = 4.920000076293945
$ws.Range("I47").Value = 17.04000091552734
$ws.Range("J47").Value = 13.27999973297119
$ws.Range("K47").Value = 23.1200008392334
$ws.Range("L47").Value = 2.680000066757202
$ws.Range("M47").Value = 36.33000183105469

# --- Replicate the date-cell format (s="2") from the existing rows onto the new ones ---
$ws.Range("A45").Copy()
$ws.Range("A46:A47").PasteSpecial(-4122)
$excel.CutCopyMode = $false


